$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("ALC")
$ws2 = $wb.Worksheets.Item("ARM")
$ws3 = $wb.Worksheets.Item("BSM")
$ws4 = $wb.Worksheets.Item("CRP")
$ws5 = $wb.Worksheets.Item("CUL")
$ws6 = $wb.Worksheets.Item("GSM")
$ws7 = $wb.Worksheets.Item("LTW")
$ws8 = $wb.Worksheets.Item("WVR")

# --- ALC ---
# row 28
$ws1.Range("H28").Value = 10899.6
$ws1.Range("I28").Value = 4121.0713
$ws1.Range("K28").Value = 4121.0713
$ws1.Range("M28").Value = -3636.0713
# row 55
$ws1.Range("H55").Value = 95.833336
$ws1.Range("I55").Value = 75.111115
$ws1.Range("J55").Value = 116.55556
$ws1.Range("K55").Value = 75.111115
$ws1.Range("L55").Value = 116.55556
$ws1.Range("M55").Value = 138.888885
$ws1.Range("N55").Value = -544.55556
# row 69
$ws1.Range("H69").Value = 8082.1177
$ws1.Range("I69").Value = 7036
$ws1.Range("K69").Value = 21108
$ws1.Range("M69").Value = -20234
# row 72
$ws1.Range("H72").Value = 8082.1177
$ws1.Range("I72").Value = 7036
$ws1.Range("K72").Value = 63324
$ws1.Range("M72").Value = -58956
# row 92
$ws1.Range("H92").Value = 1581.6364
$ws1.Range("I92").Value = 687.875
$ws1.Range("K92").Value = 687.875
$ws1.Range("M92").Value = 560.125
# row 132
$ws1.Range("H132").Value = 224765.67
$ws1.Range("I132").Value = 3270.2856
$ws1.Range("K132").Value = 9810.856800000001
$ws1.Range("M132").Value = -7280.856800000001

# --- ARM ---
# row 32
$ws2.Range("H32").Value = 3947.7334
$ws2.Range("I32").Value = 1573.8235
$ws2.Range("K32").Value = 1573.8235
$ws2.Range("M32").Value = -1286.8235
# row 61
$ws2.Range("H61").Value = 4670.25
$ws2.Range("I61").Value = 4966.7144
$ws2.Range("J61").Value = 3978.5
$ws2.Range("K61").Value = 4966.7144
$ws2.Range("L61").Value = 3978.5
$ws2.Range("M61").Value = -4754.7144
$ws2.Range("N61").Value = -4402.5
# row 63
$ws2.Range("H63").Value = 2201.9524
$ws2.Range("I63").Value = 1250.1177
$ws2.Range("J63").Value = 6247.25
$ws2.Range("K63").Value = 1250.1177
$ws2.Range("L63").Value = 6247.25
$ws2.Range("M63").Value = -564.1177
$ws2.Range("N63").Value = -7619.25
# row 66
$ws2.Range("H66").Value = 2201.9524
$ws2.Range("I66").Value = 1250.1177
$ws2.Range("J66").Value = 6247.25
$ws2.Range("K66").Value = 6250.5885
$ws2.Range("L66").Value = 31236.25
$ws2.Range("M66").Value = -2818.5885
$ws2.Range("N66").Value = -38100.25
# row 74
$ws2.Range("H74").Value = 1775.1464
$ws2.Range("I74").Value = 1461.8518
$ws2.Range("J74").Value = 2379.3572
$ws2.Range("K74").Value = 1461.8518
$ws2.Range("L74").Value = 2379.3572
$ws2.Range("M74").Value = -587.8517999999999
$ws2.Range("N74").Value = -4127.3572
# row 77
$ws2.Range("H77").Value = 1775.1464
$ws2.Range("I77").Value = 1461.8518
$ws2.Range("J77").Value = 2379.3572
$ws2.Range("K77").Value = 7309.259
$ws2.Range("L77").Value = 11896.786
$ws2.Range("M77").Value = -2941.259
$ws2.Range("N77").Value = -20632.786
# row 102
$ws2.Range("H102").Value = 1882.9412
$ws2.Range("J102").Value = 1065.6666
$ws2.Range("L102").Value = 1065.6666
$ws2.Range("N102").Value = -4309.6666
# row 122
$ws2.Range("H122").Value = 2165
$ws2.Range("I122").Value = 2165
$ws2.Range("K122").Value = 6495
$ws2.Range("M122").Value = -4045
# row 136
$ws2.Range("H136").Value = 4670.25
$ws2.Range("I136").Value = 4966.7144
$ws2.Range("J136").Value = 3978.5
$ws2.Range("K136").Value = 14900.1432
$ws2.Range("L136").Value = 11935.5
$ws2.Range("M136").Value = -12350.1432
$ws2.Range("N136").Value = -17035.5

# --- BSM ---
# row 86
$ws3.Range("H86").Value = 2709.5
$ws3.Range("I86").Value = 2420
$ws3.Range("K86").Value = 2420
$ws3.Range("M86").Value = -1297
# row 89
$ws3.Range("H89").Value = 2709.5
$ws3.Range("I89").Value = 2420
$ws3.Range("K89").Value = 12100
$ws3.Range("M89").Value = -6484
# row 105
$ws3.Range("H105").Value = 3296.2856
$ws3.Range("I105").Value = 3286.2727
$ws3.Range("K105").Value = 3286.2727
$ws3.Range("M105").Value = -1539.2727
# row 107
$ws3.Range("H107").Value = 2067.2144
$ws3.Range("I107").Value = 1091
$ws3.Range("K107").Value = 1091
$ws3.Range("M107").Value = 829

# --- CRP ---
# row 31
$ws4.Range("H31").Value = 4735.625
$ws4.Range("I31").Value = 3771.7
$ws4.Range("K31").Value = 3771.7
$ws4.Range("M31").Value = -3476.7
# row 34
$ws4.Range("H34").Value = 4735.625
$ws4.Range("I34").Value = 3771.7
$ws4.Range("K34").Value = 3771.7
$ws4.Range("M34").Value = -3569.7
# row 74
$ws4.Range("H74").Value = 28581.666
$ws4.Range("J74").Value = 28581.666
$ws4.Range("L74").Value = 28581.666
$ws4.Range("N74").Value = -30329.666
# row 77
$ws4.Range("H77").Value = 28581.666
$ws4.Range("J77").Value = 28581.666
$ws4.Range("L77").Value = 85744.99800000001
$ws4.Range("N77").Value = -94480.99800000001
# row 107
$ws4.Range("H107").Value = 753.0769
$ws4.Range("I107").Value = 639.4400000000001
$ws4.Range("J107").Value = 956
$ws4.Range("K107").Value = 639.4400000000001
$ws4.Range("L107").Value = 956
$ws4.Range("M107").Value = 1280.56
$ws4.Range("N107").Value = -4796
# row 132
$ws4.Range("H132").Value = 2718.5625
$ws4.Range("J132").Value = 2059.25
$ws4.Range("L132").Value = 6177.75
$ws4.Range("N132").Value = -11237.75
# row 134
$ws4.Range("H134").Value = 2543.9375
$ws4.Range("I134").Value = 2331.7
$ws4.Range("K134").Value = 6995.099999999999
$ws4.Range("M134").Value = -4460.099999999999

# --- CUL ---
# row 128
$ws5.Range("H128").Value = 188997
$ws5.Range("I128").Value = 188997
$ws5.Range("K128").Value = 566991
$ws5.Range("M128").Value = -562011

# --- GSM ---
# row 113
$ws6.Range("H113").Value = 933
$ws6.Range("I113").Value = 933
$ws6.Range("K113").Value = 933
$ws6.Range("M113").Value = 1237
# row 122
$ws6.Range("H122").Value = 2602.6206
$ws6.Range("I122").Value = 2845.7273
$ws6.Range("J122").Value = 2454.0557
$ws6.Range("K122").Value = 8537.1819
$ws6.Range("L122").Value = 7362.1671
$ws6.Range("M122").Value = -6087.1819
$ws6.Range("N122").Value = -12262.1671
# row 132
$ws6.Range("H132").Value = 2905.7896
$ws6.Range("I132").Value = 2983.353
$ws6.Range("K132").Value = 8950.059000000001
$ws6.Range("M132").Value = -6420.059000000001

# --- LTW ---
# row 16
$ws7.Range("H16").Value = 4095
$ws7.Range("I16").Value = 4229.357
$ws7.Range("J16").Value = 3468
$ws7.Range("K16").Value = 4229.357
$ws7.Range("L16").Value = 3468
$ws7.Range("M16").Value = -4059.357
$ws7.Range("N16").Value = -3808
# row 61
$ws7.Range("H61").Value = 5887.4546
$ws7.Range("I61").Value = 5635.263
$ws7.Range("K61").Value = 5635.263
$ws7.Range("M61").Value = -5433.263
# row 100
$ws7.Range("H100").Value = 4427.857
$ws7.Range("I100").Value = 3748.75
$ws7.Range("J100").Value = 5333.3335
$ws7.Range("K100").Value = 3748.75
$ws7.Range("L100").Value = 5333.3335
$ws7.Range("M100").Value = -3207.75
$ws7.Range("N100").Value = -6415.3335
# row 113
$ws7.Range("H113").Value = 5887.4546
$ws7.Range("I113").Value = 5635.263
$ws7.Range("K113").Value = 5635.263
$ws7.Range("M113").Value = -3465.263
# row 122
$ws7.Range("H122").Value = 4091.1667
$ws7.Range("I122").Value = 4232.8335
$ws7.Range("K122").Value = 12698.5005
$ws7.Range("M122").Value = -10248.5005
# row 132
$ws7.Range("H132").Value = 3968.6843
$ws7.Range("I132").Value = 3980.8333
$ws7.Range("J132").Value = 3947.8572
$ws7.Range("K132").Value = 11942.4999
$ws7.Range("L132").Value = 11843.5716
$ws7.Range("M132").Value = -9412.499899999999
$ws7.Range("N132").Value = -16903.5716
# row 136
$ws7.Range("H136").Value = 2644.9
$ws7.Range("I136").Value = 2750.1428
$ws7.Range("K136").Value = 8250.428400000001
$ws7.Range("M136").Value = -5700.428400000001

# --- WVR ---
# row 45
$ws8.Range("H45").Value = 17409.834
$ws8.Range("J45").Value = 18691.8
$ws8.Range("L45").Value = 18691.8
$ws8.Range("N45").Value = -19673.8
# row 100
$ws8.Range("H100").Value = 2292.818
$ws8.Range("J100").Value = 2943.8
$ws8.Range("L100").Value = 5887.6
$ws8.Range("N100").Value = -6969.6
# row 122
$ws8.Range("H122").Value = 3801.7036
$ws8.Range("I122").Value = 4138.2144
$ws8.Range("J122").Value = 3439.3076
$ws8.Range("K122").Value = 12414.6432
$ws8.Range("L122").Value = 10317.9228
$ws8.Range("M122").Value = -9964.643199999999
$ws8.Range("N122").Value = -15217.9228
